# RAD test data refresh: update the "Date" column (B) on each RAD
# worksheet with new Katalon test-run timestamps, in the order the
# test suite re-ran: Estimated, Existing, Extension, NewTaxReturn,
# Personal_EL, Personal_IND, Personal_JNT. (Personal has no data rows
# in column B, so it is left untouched.)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("B2").Value = "Mon Oct 09 23:18:02 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 23:18:42 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 23:19:23 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 09 23:20:03 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 09 23:20:43 EDT 2023"
$ws.Range("B7").Value = "Mon Oct 09 23:21:23 EDT 2023"

$ws = $wb.Worksheets.Item("Existing")
$ws.Range("B2").Value = "Mon Oct 09 23:22:03 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 23:22:42 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 23:23:22 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 09 23:24:02 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 09 23:24:42 EDT 2023"
$ws.Range("B7").Value = "Mon Oct 09 23:25:22 EDT 2023"
$ws.Range("B8").Value = "Mon Oct 09 23:26:03 EDT 2023"
$ws.Range("B9").Value = "Mon Oct 09 23:26:44 EDT 2023"
$ws.Range("B10").Value = "Mon Oct 09 23:27:25 EDT 2023"
$ws.Range("B11").Value = "Mon Oct 09 23:28:05 EDT 2023"
$ws.Range("B12").Value = "Mon Oct 09 23:28:45 EDT 2023"

$ws = $wb.Worksheets.Item("Extension")
$ws.Range("B2").Value = "Mon Oct 09 23:29:25 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 23:30:05 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 23:30:44 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 09 23:31:24 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 09 23:32:03 EDT 2023"
$ws.Range("B7").Value = "Mon Oct 09 23:32:42 EDT 2023"

$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Range("B2").Value = "Mon Oct 09 23:33:24 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 23:34:03 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 23:34:43 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 09 23:35:23 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 09 23:36:02 EDT 2023"
$ws.Range("B7").Value = "Mon Oct 09 23:36:41 EDT 2023"
$ws.Range("B8").Value = "Mon Oct 09 23:37:21 EDT 2023"
$ws.Range("B9").Value = "Mon Oct 09 23:38:01 EDT 2023"
$ws.Range("B10").Value = "Mon Oct 09 23:38:40 EDT 2023"
$ws.Range("B11").Value = "Mon Oct 09 23:39:22 EDT 2023"
$ws.Range("B12").Value = "Mon Oct 09 23:40:02 EDT 2023"
$ws.Range("B13").Value = "Mon Oct 09 23:40:42 EDT 2023"
$ws.Range("B14").Value = "Mon Oct 09 23:41:24 EDT 2023"
$ws.Range("B15").Value = "Mon Oct 09 23:42:04 EDT 2023"
$ws.Range("B16").Value = "Mon Oct 09 23:42:43 EDT 2023"

$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Range("B2").Value = "Mon Oct 09 23:43:23 EDT 2023"

$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Range("B2").Value = "Mon Oct 09 23:44:02 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 23:44:39 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 23:45:16 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 09 23:45:54 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 09 23:46:31 EDT 2023"

$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Range("B2").Value = "Mon Oct 09 23:47:11 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 23:47:54 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 23:48:39 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 09 23:49:22 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 09 23:50:06 EDT 2023"
